$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 11 (D010): Status changes from Open to Rejected, add a Coments value
$ws.Range("C11").Value = "Rejected"
$ws.Range("D11").Value = "Working as per requirement specified in assignment."

# Add a new defect row 13 (D012): Description + Status
$ws.Range("B13").Value = "On clicking FAQ, exception is occuring."
$ws.Range("C13").Value = "Open"

$ws.Range("D13").Select()
